# Project Sample Project is saved.TEST Author: admin. Type: SAVE.
# Change cell B11 on the "Rules" sheet from the text "R40" to the text "1".
#
# A plain Range.Value assignment of "1" would be auto-coerced to a number by
# Excel, but the original cell stores a shared string (t="s"), so the text
# has to be written in a way that keeps it a string. Writing it as a
# formula that evaluates to the text "1" and then collapsing the formula to
# its literal value with a values-only Paste Special preserves the text type
# (and the cell's existing style) instead of turning it into a number.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$cell = $ws.Range("B11")
$cell.Formula = '="1"'
$cell.Copy()
$cell.PasteSpecial(-4163)   # xlPasteValues
$excel.CutCopyMode = $false
